$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.376.37"
$ws.Range("E2").Value = '  +1.40%  '

$ws.Range("D3").Value = "'3.357.96"
$ws.Range("E3").Value = '  +1.22%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = "'555.41"
$ws.Range("E5").Value = '  +1.17%  '

$ws.Range("D6").Value = "'173.71"
$ws.Range("E6").Value = '  +0.40%  '

$ws.Range("D7").Value = "'0.619"
$ws.Range("E7").Value = '  +1.45%  '

$ws.Range("D8").Value = "'3.348.12"
$ws.Range("E8").Value = '  +1.13%  '

$ws.Range("E9").Value = '  +0.08%  '

$ws.Range("B10").Value = 'Cardano'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D10").Value = "'0.627"
$ws.Range("E10").Value = '  +2.53%  '

$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").Value = "'0.163"
$ws.Range("E11").Value = '  +7.27%  '

$ws.Range("D12").Value = "'53.59"
$ws.Range("E12").Value = '  +0.80%  '

$ws.Range("E13").Value = '  +2.98%  '

$ws.Range("D14").Value = "'9.03"
$ws.Range("E14").Value = '  +1.50%  '

$ws.Range("D15").Value = "'3.897.63"
$ws.Range("E15").Value = '  +1.39%  '

$ws.Range("D16").Value = "'18.16"
$ws.Range("E16").Value = '  -0.11%  '

$ws.Range("E17").Value = '  +1.49%  '

$ws.Range("D18").Value = "'3.358.31"
$ws.Range("E18").Value = '  +1.36%  '

$ws.Range("D19").Value = "'64.446.13"
$ws.Range("E19").Value = '  +1.61%  '

$ws.Range("E20").Value = '  -0.03%  '

$ws.Range("E21").Value = '  +1.72%  '

$ws.Range("D22").Value = "'457.27"
$ws.Range("E22").Value = '  +7.84%  '

$ws.Range("E23").Value = '  +10.20%  '

$ws.Range("D24").Value = "'4.08"
$ws.Range("E24").Value = '  +0.68%  '

$ws.Range("D25").Value = "'85.83"
$ws.Range("E25").Value = '  +3.32%  '

$ws.Range("D26").Value = "'13.55"
$ws.Range("E26").Value = '  +1.36%  '

$ws.Range("E27").Value = '  +8.41%  '

$ws.Range("D28").Value = "'10.72"
$ws.Range("E28").Value = '  +0.96%  '

$ws.Range("E29").Value = '  +0.22%  '

$ws.Range("D30").Value = "'30.38"
$ws.Range("E30").Value = '  +4.16%  '

$ws.Range("D31").Value = "'6.66"
$ws.Range("E31").Value = '  +3.18%  '

$ws.Range("E32").Value = '  +0.36%  '

$ws.Range("D33").Value = "'571.33"
$ws.Range("E33").Value = '  -1.00%  '

$ws.Range("D34").Value = "'60.92"
$ws.Range("E34").Value = '  +4.70%  '

$ws.Range("E35").Value = '  +1.08%  '

$ws.Range("E36").Value = '  -0.03%  '

$ws.Range("D37").Value = "'3.63"
$ws.Range("E37").Value = '  +4.56%  '

$ws.Range("E38").Value = '  -4.90%  '

$ws.Range("D39").Value = "'35.22"
$ws.Range("E39").Value = '  +0.33%  '

$ws.Range("D40").Value = "'0.0₃0737"
$ws.Range("E40").Value = '  -0.26%  '

$ws.Range("E41").Value = '  +0.75%  '

$ws.Range("E42").Value = '  +0.09%  '

$ws.Range("D43").Value = "'3.068.10"
$ws.Range("E43").Value = '  -1.75%  '

$ws.Range("E44").Value = '  +0.63%  '

$ws.Range("D45").Value = "'0.0412"
$ws.Range("E45").Value = '  +2.69%  '

$ws.Range("E46").Value = '  +4.15%  '

$ws.Range("D47").Value = "'2.43"
$ws.Range("E47").Value = '  +0.58%  '

$ws.Range("D48").Value = "'3.14"
$ws.Range("E48").Value = '  -0.44%  '

$ws.Range("D50").Value = "'138.11"
$ws.Range("E50").Value = '  +2.33%  '

$ws.Range("D51").Value = "'8.11"
$ws.Range("E51").Value = '  +1.05%  '
